$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 3006
$ws1.Range("F10").Value = 7080
$ws1.Range("F14").Value = 633
$ws1.Range("F16").Value = 2287
$ws1.Range("F17").Value = 1548
$ws1.Range("F20").Value = 160
$ws1.Range("F22").Value = 209
$ws1.Range("F24").Value = 57
$ws1.Range("F25").Value = 57
$ws1.Range("F26").Value = 1804
$ws1.Range("F30").Value = 1695
$ws1.Range("F41").Value = 43

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F13").Value = 73

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F6").Value = 1736
$ws3.Range("F8").Value = 2797
$ws3.Range("F9").Value = 1063
$ws3.Range("F10").Value = 976
$ws3.Range("F12").Value = 339
$ws3.Range("F13").Value = 1658
$ws3.Range("F14").Value = 7627

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 3006
$ws4.Range("F6").Value = 1736
$ws4.Range("F7").Value = 2797
$ws4.Range("F8").Value = 7080
$ws4.Range("F9").Value = 1063
$ws4.Range("F12").Value = 633
$ws4.Range("F14").Value = 2287
$ws4.Range("F15").Value = 1548
$ws4.Range("F18").Value = 160
$ws4.Range("F21").Value = 57
$ws4.Range("F22").Value = 57
$ws4.Range("F23").Value = 1804
$ws4.Range("F24").Value = 73
$ws4.Range("F27").Value = 1695
$ws4.Range("F41").Value = 43
